# Organised Full Item list
# Fill in the "Attument" (attunement) column (F) on Sheet1 with "N" for every
# row whose value is currently blank, and correct the two rows that were
# mistakenly marked "Y" (rows 17 and 18) to "N".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(5,6,7,9,10,12,13,16,17,18,19,20,21,22,23,24,25,26,29,30,31,32,33,34,37,39,40,43,44,45,47,50,51,52,53,55,56,57,58,61,62,64,69,73,74,76,80,87,91,92,93,95,96,97,101)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "N"
}
